$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("IAM")

# Append the x-1p-zuul-host header to the X-1P-User rows that had their
# test requests failing (OPQA-3768 related resetpassword scripts).
$zuulSuffix = "||x-1p-zuul-host=http://localhost:7001"

$rowsToAppend = @(34, 35, 37, 38, 39, 40, 43, 44, 45, 46)
foreach ($r in $rowsToAppend) {
    $cell = $ws.Cells.Item($r, 6)  # column F
    $cell.Value = $cell.Value2 + $zuulSuffix
}

# Fix the expected validation for the short-term token test (row 62, col J)
$ws.Cells.Item(62, 10).Value = "status=200"

# Update the selection to reflect the last edited range
$ws.Range("L2:L89").Select()
